$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bill Summary")

# --- Column C: "Qty executed upto date" (plain numeric cells) ---
$ws.Range("C8").Value = 29
$ws.Range("C9").Value = 70
$ws.Range("C10").Value = 79
$ws.Range("C11").Value = 5
$ws.Range("C12").Value = 34
$ws.Range("C13").Value = 52
$ws.Range("C14").Value = 39
$ws.Range("C15").Value = 85
$ws.Range("C16").Value = 10
$ws.Range("C17").Value = 63

# --- Columns G & H: "Upto date Amount" / grand-total cells store
# numeric-looking values as TEXT. Pre-format each cell as text so Excel
# does not silently coerce the assigned strings back into numbers, then
# clear the formatting afterwards so no stray cell style is left behind
# (the original workbook has no cell-level styles at all). ---
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "17920.00"
$ws.Range("G9").ClearFormats()

$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "37288.00"
$ws.Range("G10").ClearFormats()

$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "3310.00"
$ws.Range("G11").ClearFormats()

$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "7072.00"
$ws.Range("G13").ClearFormats()

$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "897.00"
$ws.Range("G14").ClearFormats()

$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "66487.00"
$ws.Range("G19").ClearFormats()

$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "66487.00"
$ws.Range("H19").ClearFormats()

$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "66487.00"
$ws.Range("G21").ClearFormats()

$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "66487.00"
$ws.Range("H21").ClearFormats()
